# Update the Scoring column (C) on the "Scores" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

$ws.Range("C2").Value = "Neutral"
$ws.Range("C3").Value = "Very Grandiose"
$ws.Range("C4").Value = "Neutral"
$ws.Range("C5").Value = "Neutral"
$ws.Range("C6").Value = "Neutral"
$ws.Range("C7").Value = "Neutral"
$ws.Range("C8").Value = "Neutral"
$ws.Range("C9").Value = "Neutral"
$ws.Range("C10").Value = "Neutral"
$ws.Range("C11").Value = "Neutral"
$ws.Range("C12").Value = "Neutral"
$ws.Range("C13").Value = "Neutral"
$ws.Range("C14").Value = "Neutral"
$ws.Range("C15").Value = "Neutral"
$ws.Range("C16").Value = "Neutral"
$ws.Range("C17").Value = "Neutral"
$ws.Range("C18").Value = "Neutral"
$ws.Range("C19").Value = "Neutral"
$ws.Range("C20").Value = "Neutral"
